$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Find the next empty row in column A (appends after the last existing log entry)
$lastRow = $ws.Cells.Item($ws.Rows.Count, 1).End(-4162).Row
$nextRow = $lastRow + 1

# Append the new timestamp log entry as plain text
$ws.Cells.Item($nextRow, 1).Value = "2025-07-28 10:44:00"
